$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# 1) Insert a new column B ("Week_Start_Date") - shifts ASIN..is_holiday_week from B..I to C..J
$ws.Columns("B").Insert()

# 2) Header for the new column
$ws.Range("B1").Value = "Week_Start_Date"

# 3) Week_Start_Date values for each of the 16 data rows (stored as plain text, matching "YYYY-MM-DD")
$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("B$row")
    # Force text format first so the ISO date-like string isn't auto-converted to a date serial
    $cell.NumberFormat = "@"
    $cell.Value = $weekStartDates[$i]
    # Restore the default "Normal" style so no stray number-format style sticks to the cell
    $cell.Style = "Normal"
}

# 4) Normalize the "Week" labels from zero-padded (W01..W09) to unpadded (W1..W9)
for ($i = 1; $i -le 9; $i++) {
    $row = $i + 1
    $ws.Range("A$row").Value = "W$i"
}

# 5) Fix MyForecast (now column D) value for week 4 (row 5): 4 -> 3
$ws.Range("D5").Value = 3

# 6) is_holiday_week (now column J) must be a boolean, not a number
for ($row = 2; $row -le 17; $row++) {
    $ws.Range("J$row").Value = $false
}

# 7) Update the Summary sheet totals to reflect the corrected forecast
# (these cells are stored as text in the workbook, so force text format to
#  avoid the numeric-looking strings being auto-converted to numbers)
$summary = $wb.Worksheets.Item("Summary")
$summaryUpdates = @{ "B9" = "72"; "B10" = "29"; "B11" = "10" }
foreach ($addr in $summaryUpdates.Keys) {
    $cell = $summary.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $summaryUpdates[$addr]
    $cell.Style = "Normal"
}
